$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Update the "ParticipantsTab" Cypher query text (B2) ---
$newQuery = @'
MATCH (p:participant)-->(s:study)
OPTIONAL MATCH (samp:sample)-->(p)
OPTIONAL MATCH (p)<--(diag:diagnosis)
OPTIONAL MATCH (samp)<--(f:file)
OPTIONAL MATCH (f)<--(g:genomic_info)
WITH s, p, samp, f, g, diag
WHERE g.library_source in['Bulk Whole Cell']
WITH p
OPTIONAL MATCH (p)-->(s:study)
OPTIONAL MATCH (samp:sample)-->(p)
WITH s, p, apoc.coll.sort(collect(distinct samp.sample_id)) as samp
RETURN 
coalesce(p.participant_id,'') as `Participant ID`,
coalesce(s.study_name, '') as `Study Name`,
coalesce(s.phs_accession,'') as `Accession`,
coalesce(p.gender,'') as `Gender`,
coalesce(apoc.text.join(samp, ','), '') as `Samples`
ORDER BY p.participant_id
LIMIT 100
'@

$ws.Range("B2").Value = $newQuery

# --- Bump the base font size used on the sheet (12 -> 15 pt) ---
# Apply only to the cells that already carry content/formatting so we
# don't materialize brand-new blank cells that weren't in the sheet.
$ws.Range("A1:E1").Font.Size = 15
$ws.Range("A2").Font.Size = 15
$ws.Range("D2:E2").Font.Size = 15
$ws.Range("A3").Font.Size = 15
$ws.Range("D3:E3").Font.Size = 15
$ws.Range("A4").Font.Size = 15
$ws.Range("D4:E4").Font.Size = 15
$ws.Range("B2:C4").Font.Size = 15
$ws.Range("B5:C5").Font.Size = 15
$ws.Range("C6").Font.Size = 15

# --- Row heights grow to fit the bigger font / longer wrapped query text ---
$ws.Rows.Item(2).RowHeight = 390
$ws.Rows.Item(3).RowHeight = 292.5
$ws.Rows.Item(4).RowHeight = 292.5

# --- Selection moved to H3 ---
$null = $ws.Range("H3").Select()
